# Trade #80 (MarketMaking) closes with an early_exit at 2026-02-17 21:08:30,
# and a new trade (#113) is opened right after. Propagate both the trade-log
# rows and the summary/strategy roll-up numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: overall trade count + win rate shift with the new trade
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 80      # Total Trades
$summary.Range("B9").Value = 47.5    # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 5) trade count + win rate
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 47       # Trades
$status.Range("G5").Value = 51.06    # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet: close trade #80 (row 81) and append new trade #113
# (row 114). Columns: A#, B Date, C Time, D Strategy, E Side, F Entry,
# G Exit, H Status, I P&L%, J P&L$, K Capital After, L Exit Reason,
# M Duration, N Entry Slip, O Exit Slip, P Confidence, Q Entry Reason
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G81").Value = 0.02
$allTrades.Range("H81").Value = "CLOSED"
$allTrades.Range("K81").Value = 100.61
$allTrades.Range("L81").Value = "early_exit"
$allTrades.Range("M81").Value = 0.19

# New row 114 - copy the date cell from an existing same-day row so it
# stays a plain text value instead of being reinterpreted as a date.
$allTrades.Range("B113").Copy($allTrades.Range("B114"))
$allTrades.Cells.Item(114, 1).Value = 113
$allTrades.Cells.Item(114, 3).Value = "21:08:24"
$allTrades.Cells.Item(114, 4).Value = "MarketMaking"
$allTrades.Cells.Item(114, 5).Value = "DOWN"
$allTrades.Cells.Item(114, 6).Value = 0.02
$allTrades.Cells.Item(114, 8).Value = "OPEN"
$allTrades.Cells.Item(114, 9).Value = 0
$allTrades.Cells.Item(114, 10).Value = 0
$allTrades.Cells.Item(114, 11).Value = 100.6114872031006
$allTrades.Cells.Item(114, 13).Value = 0
$allTrades.Cells.Item(114, 14).Value = 0
$allTrades.Cells.Item(114, 15).Value = 0
$allTrades.Cells.Item(114, 16).Value = 0.6
$allTrades.Cells.Item(114, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet: close trade #80 (row 48) and append new trade #113
# (row 81). Columns: A#, B Date, C Time, D Strategy, E Side, F Entry,
# G Exit, H Status, I P&L%, J P&L$, K Capital After, L Entry Slip,
# M Exit Slip, N Confidence, O Entry Reason, P Exit Reason, Q Duration
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Range("G48").Value = 0.02
$marketMaking.Range("H48").Value = "CLOSED"
$marketMaking.Range("K48").Value = 100.61
$marketMaking.Range("P48").Value = "early_exit"
$marketMaking.Range("Q48").Value = 0.19

# New row 81 - same date-copy trick as above.
$marketMaking.Range("B80").Copy($marketMaking.Range("B81"))
$marketMaking.Cells.Item(81, 1).Value = 113
$marketMaking.Cells.Item(81, 3).Value = "21:08:24"
$marketMaking.Cells.Item(81, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(81, 5).Value = "DOWN"
$marketMaking.Cells.Item(81, 6).Value = 0.02
$marketMaking.Cells.Item(81, 8).Value = "OPEN"
$marketMaking.Cells.Item(81, 9).Value = 0
$marketMaking.Cells.Item(81, 10).Value = 0
$marketMaking.Cells.Item(81, 11).Value = 100.6114872031006
$marketMaking.Cells.Item(81, 12).Value = 0
$marketMaking.Cells.Item(81, 13).Value = 0
$marketMaking.Cells.Item(81, 14).Value = 0.6
$marketMaking.Cells.Item(81, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(81, 17).Value = 0
